$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.347.15"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "1.822.98"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'313.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4471"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.32%  "

$ws.Range("D8").Value = "'0.3755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.16%  "

$ws.Range("D9").Value = "'0.07521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.49%  "

$ws.Range("D10").Value = "'0.8856"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.72%  "

$ws.Range("D11").Value = "'21.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.69%  "

$ws.Range("D12").Value = "1.828.02"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "'6.760"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.59%  "

$ws.Range("D14").Value = "'93.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.00%  "

$ws.Range("D15").Value = "'5.403"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.03%  "

$ws.Range("D16").Value = "'0.07108"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "'0.000008811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").Value = "'15.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.83%  "

$ws.Range("D21").Value = "27.352.64"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").Value = "'5.257"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("D23").Value = "'10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").Value = "2.055.46"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "'1.966"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("D26").Value = "'2.374"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.96%  "

$ws.Range("D27").Value = "'151.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'18.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.49%  "

$ws.Range("D29").Value = "'5.361"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'118.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("D31").Value = "'0.08815"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").Value = "'0.7900"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.67%  "

$ws.Range("D33").Value = "'1.197"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "

$ws.Range("D34").Value = "'4.513"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "'1.112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "

$ws.Range("D38").Value = "'0.01995"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.48%  "

$ws.Range("D39").Value = "'0.05332"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.95%  "

$ws.Range("D40").Value = "'7.387"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.15%  "

$ws.Range("D41").Value = "'0.5314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.32%  "

$ws.Range("D42").Value = "'0.1726"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("D43").Value = "'2.859"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").Value = "'2.323"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +19.76%  "

$ws.Range("E45").Value = "  +1.99%  "

$ws.Range("D46").Value = "'0.5122"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.26%  "

$ws.Range("D47").Value = "'10.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("D48").Value = "'105.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").Value = "'1.701"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").Value = "'0.06380"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
